# edit.ps1
# Applies the STEP 1B acceptance freeze changes to l21_rerun_results.xlsx
# - Fix CAP-001 metadata (role -> OBSERVE, mutability -> READ, bound_surface -> L21-EVD-R)
# - Re-run full pipeline output re-ordering / value corrections per domain block

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "CAP-011"
$ws.Range("D2").Value = "BOUNDED"
$ws.Range("A3").Value = "CAP-016"
$ws.Range("C3").Value = "CONTROL"
$ws.Range("A4").Value = "CAP-002"
$ws.Range("C4").Value = "CONTROL"
$ws.Range("A5").Value = "CAP-008"
$ws.Range("E5").Value = "READ"
$ws.Range("F5").Value = "L21-ACT-R"
$ws.Range("A6").Value = "CAP-012"
$ws.Range("D6").Value = "STRICT"
$ws.Range("A7").Value = "CAP-020"
$ws.Range("C7").Value = "CONTROL"
$ws.Range("D7").Value = "BOUNDED"
$ws.Range("A8").Value = "CAP-021"
$ws.Range("C8").Value = "CONTROL"
$ws.Range("A9").Value = "CAP-001"
$ws.Range("C9").Value = "OBSERVE"
$ws.Range("D9").Value = "STRICT"
$ws.Range("F9").Value = "L21-EVD-R"
$ws.Range("C10").Value = "OBSERVE"
$ws.Range("E10").Value = "GOVERN"
$ws.Range("F10").Value = "L21-EVD-R"
$ws.Range("C11").Value = "CONTROL"
$ws.Range("C12").Value = "OBSERVE"
$ws.Range("F12").Value = "L21-EVD-R"
$ws.Range("C13").Value = "CONTROL"
$ws.Range("C14").Value = "CONTROL"
$ws.Range("A16").Value = "CAP-002"
$ws.Range("C16").Value = "CONTROL"
$ws.Range("D16").Value = "BOUNDED"
$ws.Range("A17").Value = "CAP-005"
$ws.Range("C17").Value = "CONTROL"
$ws.Range("E17").Value = "GOVERN"
$ws.Range("A18").Value = "CAP-009"
$ws.Range("A19").Value = "CAP-021"
$ws.Range("C19").Value = "CONTROL"
$ws.Range("E19").Value = "READ"
$ws.Range("F19").Value = "L21-ACT-R"
$ws.Range("A20").Value = "CAP-001"
$ws.Range("C20").Value = "OBSERVE"
$ws.Range("D20").Value = "STRICT"
$ws.Range("F20").Value = "L21-EVD-R"
$ws.Range("A21").Value = "CAP-002"
$ws.Range("C21").Value = "CONTROL"
$ws.Range("D21").Value = "BOUNDED"
$ws.Range("A22").Value = "CAP-021"
$ws.Range("C22").Value = "CONTROL"
$ws.Range("E22").Value = "READ"
$ws.Range("F22").Value = "L21-ACT-R"
$ws.Range("A23").Value = "CAP-001"
$ws.Range("C23").Value = "OBSERVE"
$ws.Range("D23").Value = "STRICT"
$ws.Range("F23").Value = "L21-EVD-R"
$ws.Range("C24").Value = "CONTROL"
$ws.Range("A28").Value = "CAP-005"
$ws.Range("C28").Value = "CONTROL"
$ws.Range("D28").Value = "BOUNDED"
$ws.Range("E28").Value = "GOVERN"
$ws.Range("A29").Value = "CAP-011"
$ws.Range("A30").Value = "CAP-003"
$ws.Range("C30").Value = "OBSERVE"
$ws.Range("E30").Value = "READ"
$ws.Range("F30").Value = "L21-EVD-R"
$ws.Range("A31").Value = "CAP-016"
$ws.Range("C31").Value = "CONTROL"
$ws.Range("E31").Value = "WRITE"
$ws.Range("F31").Value = "L21-ACT-W"
$ws.Range("A32").Value = "CAP-001"
$ws.Range("C32").Value = "OBSERVE"
$ws.Range("D32").Value = "STRICT"
$ws.Range("E32").Value = "READ"
$ws.Range("F32").Value = "L21-EVD-R"
$ws.Range("C34").Value = "OBSERVE"
$ws.Range("E34").Value = "GOVERN"
$ws.Range("F34").Value = "L21-EVD-R"
